$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D, shifting D:K to E:L
$ws.Columns("D").Insert()

# Copy number formats from the (now shifted) E column into the new D column
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column D with its own values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1547100
$ws.Range("D9").Value = 1978400
$ws.Range("D10").Value = -431300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 2329900
$ws.Range("D18").Value = -782800
$ws.Range("D20").Value = 1389500
$ws.Range("D21").Value = 938200
$ws.Range("D22").Value = 41000
$ws.Range("D23").Value = 565700
$ws.Range("D24").Value = 107700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 458000
$ws.Range("D27").Value = 458000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1389500
$ws.Range("D33").Value = 458000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 458000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 344400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 49222300
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 291200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 61625600
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 57986200
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 737600
$ws.Range("D62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 59226500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1550000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2399100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 458000
$ws.Range("D83").Value = 331500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 43200
$ws.Range("D91").Value = -4300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2408300
$ws.Range("D96").Value = -25300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 1275500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -1089600

# Two cells in row 89 and row 94 received revised historical values
$ws.Range("E89").Value = 1923800
$ws.Range("F89").Value = 1416400
$ws.Range("E94").Value = -2593400
$ws.Range("F94").Value = -4501100
